$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column I (9th column) header was blank; set it to "Debit" to match the
# new "Debit" column inserted between "Transaction Date" and "Credit".
$ws.Range("I1").Value = "Debit"
